$wb = $excel.ActiveWorkbook

# --- "Birth outcomes & risks" sheet: add two new rows for OR wasting ---
$ws = $wb.Worksheets.Item("Birth outcomes & risks")

# Insert two blank rows above the current row 5 (pushes old rows 5-15 down to 7-17)
$ws.Rows("5:6").Insert()

# New row 5: OR wasting (high)
$ws.Range("B5").Value = "OR wasting (high)"
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 1

# New row 6: OR wasting (moderate)
$ws.Range("B6").Value = "OR wasting (moderate)"
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 1

# Style the new data cells: Calibri 11, yellow fill, right-aligned (matches the
# sheet's other "header value" rows, e.g. row 4)
$newRange = $ws.Range("C5:F6")
$newRange.Font.Name = "Calibri"
$newRange.Font.Size = 11
$newRange.Interior.Color = 65535
$newRange.HorizontalAlignment = -4152

# --- View-state bookkeeping to match the saved workbook state ---

# "Distributions" sheet selection
$wsDist = $wb.Worksheets.Item("Distributions")
[void]$wsDist.Range("G8:G11").Select()

# "Interventions for children" sheet selection (was the active tab before the edit)
$wsInt = $wb.Worksheets.Item("Interventions for children")
[void]$wsInt.Range("A6").Select()

# "Birth outcomes & risks" becomes the active tab, with new selection
[void]$ws.Range("G6").Select()
